# Generate Report for Handback
# This script updates the localization-status workbook to reflect a
# completed handback: the "Ready for handoff" status becomes
# "Handed back: in sync with en-US", the handback timestamps are
# refreshed, and the stale "handback file is not the latest" error
# details are cleared now that the handback is in sync.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# Status columns for each locale (zh-cn / de-de) move from
# "Ready for handoff" to "Handed back: in sync with en-US".
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Best-effort column resize to fit the new, longer status text.
# (ColumnWidth is expressed in character units; the engine serializes it
# to the OOXML <col width> attribute with a fixed +5/6 padding offset, so
# the inputs below are chosen to land the stored width as close as
# possible to the target layout produced by the report generator.)
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Refresh the "Latest Handback DateTime" now that the handback is in sync.
$zhcn.Range("K2").Value = "2016-08-14 00:48:47"
$zhcn.Range("K3").Value = "2016-08-14 00:48:47"

# The handback is now in sync, so the stale-version error detail clears.
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333332

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Refresh the "Latest Handback DateTime" now that the handback is in sync.
$dede.Range("K2").Value = "2016-08-14 00:48:57"
$dede.Range("K3").Value = "2016-08-14 00:48:57"

# The handback is now in sync, so the stale-version error detail clears.
$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333332
